$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: B1 text changes, and the header area is widened out to
#     column N (blank, formatted cells C1:N1 matching B1's fill/format) ---
$ws.Range("B1").Value = "123"
$ws.Range("B1").Copy()
$ws.Range("C1:N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1:N1").ClearContents()
$excel.CutCopyMode = 0

# --- Row 5: B5 becomes a text "1" (shared string) instead of a number,
#     and C5 (previously 5) is removed entirely ---
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1"
$ws.Range("C5").Clear()

# --- Column widths for the newly-visible columns B:N ---
$ws.Columns.Item(2).ColumnWidth = 10.6328125
$ws.Columns.Item(3).ColumnWidth = 10.453125
$ws.Columns.Item(4).ColumnWidth = 11.54296875
$ws.Columns.Item(5).ColumnWidth = 11.54296875
$ws.Columns.Item(6).ColumnWidth = 11.54296875
$ws.Columns.Item(7).ColumnWidth = 11.90625
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 11.7265625
$ws.Columns.Item(10).ColumnWidth = 12
$ws.Columns.Item(11).ColumnWidth = 11.6328125
$ws.Columns.Item(12).ColumnWidth = 12.26953125
$ws.Columns.Item(13).ColumnWidth = 11.81640625
$ws.Columns.Item(14).ColumnWidth = 12.36328125

# --- View state: zoom to 85% and move the selection to E4 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("E4").Select()
